# aggiornamento fino a 9 agosto 2021
# Extend the daily log table on Sheet1 with new rows, continuing the
# existing date series (column A, serial dates) by one day per row and
# carrying the same zero values / formatting used for the prior rows
# (columns B, C, D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row/column of the existing table.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Last day currently present (column A) and the date it should end on.
$lastSerial = $ws.Cells.Item($lastRow, 1).Value2
$targetSerial = 44417   # 2021-08-09

$sourceRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol))

for ($serial = $lastSerial + 1; $serial -le $targetSerial; $serial++) {
    $r = $lastRow + ($serial - $lastSerial)
    $destRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol))

    # Copy the last row (values + formatting) as the template for the new row...
    $sourceRange.Copy($destRange) | Out-Null

    # ...then overwrite the date cell with the next day's serial value.
    $ws.Cells.Item($r, 1).Value2 = $serial
}
